# Generate Report for Archive
#
# The handback/status report moved from "Ready for handoff" to
# "In Translation" for both tracked files. Update the shared
# "Status" text everywhere it is used:
#   - Overview sheet: per-locale status columns (zh-cn, de-de)
#   - zh-cn sheet: Status column
#   - de-de sheet: Status column
# Shrinking the text also shrinks the (auto-fitted) column widths for
# those status columns, so re-apply the column width after the edit.

$wb = $excel.ActiveWorkbook

$old = "Ready for handoff"
$new = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: columns E (zh-cn) and F (de-de) hold the per-locale status
$wsOverview.Range("E2").Value = $new
$wsOverview.Range("F2").Value = $new
$wsOverview.Range("E3").Value = $new
$wsOverview.Range("F3").Value = $new

# zh-cn sheet: column C is "Status"
$wsZhCn.Range("C2").Value = $new
$wsZhCn.Range("C3").Value = $new

# de-de sheet: column C is "Status"
$wsDeDe.Range("C2").Value = $new
$wsDeDe.Range("C3").Value = $new

# The status columns were sized to fit their (now shorter) content.
$wsOverview.Columns("E:F").ColumnWidth = 12.5
$wsZhCn.Columns("C:C").ColumnWidth = 12.5
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
